$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2024-12-05 Thursday" "2024-12-06 Friday"

Replace-Text "918÷4=229, 2" "752÷4=188, 0"
Replace-Text "374÷3=124, 2" "147÷8=18, 3"
Replace-Text "221÷9=24, 5" "628÷5=125, 3"
Replace-Text "793÷6=132, 1" "776÷8=97, 0"
Replace-Text "367÷2=183, 1" "820÷9=91, 1"

Replace-Text "229÷3=76, 1" "554÷8=69, 2"
Replace-Text "600÷6=100, 0" "401÷4=100, 1"
Replace-Text "640÷6=106, 4" "619÷9=68, 7"
Replace-Text "560÷6=93, 2" "818÷3=272, 2"
Replace-Text "756÷8=94, 4" "258÷5=51, 3"

Replace-Text "359÷9=39, 8" "687÷2=343, 1"
Replace-Text "830÷3=276, 2" "694÷7=99, 1"
Replace-Text "438÷3=146, 0" "166÷4=41, 2"
Replace-Text "670÷3=223, 1" "605÷3=201, 2"
Replace-Text "904÷7=129, 1" "626÷5=125, 1"

Replace-Text "679÷5=135, 4" "974÷4=243, 2"
Replace-Text "352÷2=176, 0" "221÷3=73, 2"
Replace-Text "559÷5=111, 4" "182÷5=36, 2"
Replace-Text "635÷3=211, 2" "974÷7=139, 1"
Replace-Text "938÷8=117, 2" "702÷2=351, 0"

Replace-Text "289÷7=41, 2" "826÷7=118, 0"
Replace-Text "502÷2=251, 0" "532÷3=177, 1"
Replace-Text "822÷9=91, 3" "872÷5=174, 2"
Replace-Text "584÷3=194, 2" "587÷6=97, 5"
Replace-Text "988÷4=247, 0" "289÷9=32, 1"
